# Add a new "Before" worksheet after the existing "CellListener" sheet.
# This is the template's "before" example used to show the raw/unresolved
# cell content that CellListener will later replace.
$wb = $excel.ActiveWorkbook
$cellListenerSheet = $wb.Worksheets.Item(1)

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $cellListenerSheet)
$ws.Name = "Before"

# A1: label; B1: placeholder text that CellListener will overwrite;
# B2: explanatory note (kept as literal text, not a live expression).
$ws.Range("A1").Value = "State Name:"
$ws.Range("B1").Value = "Anything here; CellListener will replace!"
$ws.Range("B2").Value = 'The CellListener will replace the above content with ${california.name}'

$ws.Columns.Item(1).ColumnWidth = 17.6
$ws.Columns.Item(2).ColumnWidth = 35.9

# Keep "CellListener" as the active/selected sheet, as it was originally.
$cellListenerSheet.Activate()
